$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13, pushing existing rows 13-62 down to 14-63.
# Excel inherits the formatting (e.g. the date style on column D) from
# the row above, matching the surrounding rows.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with this week's new data point
# (the rest of the columns mirror the constant values used by every
# other row in this data set).
$ws.Range("A13").Value = 6
$ws.Range("B13").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 45037
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100112035
$ws.Range("G13").Value = "Bruselas (repollito)"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 20000
$ws.Range("L13").Value = 22000
$ws.Range("M13").Value = 20850
$ws.Range("N13").Value = "$/malla 15 kilos"
$ws.Range("O13").Value = "Provincia de Quillota"
$ws.Range("P13").Value = 1390
$ws.Range("Q13").Value = 15
$ws.Range("R13").Value = "Hortaliza"
